# Apply the "Updated symbol list" commit changes to the crypto price sheet.
# Every data cell in columns B:G is stored as text (inlineStr) in the workbook,
# so we prefix each written value with a leading apostrophe to force Excel to
# keep storing it as text instead of auto-coercing numeric-looking strings
# (prices, the hour counter, etc.) into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.15"
$ws.Range("G2").Value = "'9"
# Row 3
$ws.Range("G3").Value = "'9"
# Row 4
$ws.Range("D4").Value = "'5.386"
$ws.Range("G4").Value = "'9"
# Row 5
$ws.Range("D5").Value = "'0.05979"
$ws.Range("G5").Value = "'9"
# Row 6
$ws.Range("D6").Value = "'3.393"
$ws.Range("G6").Value = "'9"
# Row 7
$ws.Range("G7").Value = "'9"
# Row 8
$ws.Range("D8").Value = "'0.8110"
$ws.Range("G8").Value = "'9"
# Row 9
$ws.Range("D9").Value = "'0.9627"
$ws.Range("G9").Value = "'9"
# Row 10
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1430"
$ws.Range("E10").Value = "'9WazirXWRX"
$ws.Range("G10").Value = "'9"
# Row 11
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07414"
$ws.Range("E11").Value = "'10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "'9"
# Row 12
$ws.Range("B12").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03397"
$ws.Range("E12").Value = "'11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "'9"
# Row 13
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03058"
$ws.Range("E13").Value = "'12BitrueCoinBTR"
$ws.Range("G13").Value = "'9"
# Row 14
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09416"
$ws.Range("E14").Value = "'13BitMartTokenBMX"
$ws.Range("G14").Value = "'9"
# Row 15
$ws.Range("B15").Value = "'MCDex"
$ws.Range("C15").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'4.000"
$ws.Range("E15").Value = "'14MCDexMCB"
$ws.Range("G15").Value = "'9"
# Row 16
$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001597"
$ws.Range("E16").Value = "'15BitForexTokenBF"
$ws.Range("G16").Value = "'9"
# Row 17
$ws.Range("B17").Value = "'CoinExToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04808"
$ws.Range("E17").Value = "'16CoinExTokenCET"
$ws.Range("G17").Value = "'9"
# Row 18
$ws.Range("B18").Value = "'One"
$ws.Range("C18").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005882"
$ws.Range("E18").Value = "'17OneONE"
$ws.Range("G18").Value = "'9"
# Row 19
$ws.Range("D19").Value = "'0.006217"
$ws.Range("G19").Value = "'9"
# Row 20
$ws.Range("D20").Value = "'0.005084"
$ws.Range("G20").Value = "'9"
# Row 21
$ws.Range("D21").Value = "'0.0009848"
$ws.Range("G21").Value = "'9"
# Row 22
$ws.Range("D22").Value = "'0.0001001"
$ws.Range("G22").Value = "'9"
# Row 23
$ws.Range("D23").Value = "'3.749"
$ws.Range("G23").Value = "'9"
# Row 24
$ws.Range("G24").Value = "'9"
# Row 25
$ws.Range("G25").Value = "'9"
# Row 26
$ws.Range("G26").Value = "'9"
# Row 27
$ws.Range("G27").Value = "'9"
# Row 28
$ws.Range("G28").Value = "'9"
# Row 29
$ws.Range("G29").Value = "'9"
# Row 30
$ws.Range("G30").Value = "'9"
# Row 31
$ws.Range("G31").Value = "'9"
# Row 32
$ws.Range("G32").Value = "'9"
# Row 33
$ws.Range("G33").Value = "'9"
# Row 34
$ws.Range("G34").Value = "'9"
# Row 35
$ws.Range("G35").Value = "'9"
# Row 36
$ws.Range("G36").Value = "'9"
# Row 37
$ws.Range("G37").Value = "'9"
# Row 38
$ws.Range("G38").Value = "'9"
# Row 39
$ws.Range("G39").Value = "'9"
# Row 40
$ws.Range("D40").Value = "'0.03960"
$ws.Range("G40").Value = "'9"
# Row 41
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1076"
$ws.Range("E41").Value = "'40BKEXTokenBKK"
$ws.Range("G41").Value = "'9"
# Row 42
$ws.Range("B42").Value = "'CEJI"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002722"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("G42").Value = "'9"
# Row 43
$ws.Range("B43").Value = "'KickToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003017"
$ws.Range("E43").Value = "'42KickTokenKICK"
$ws.Range("G43").Value = "'9"
# Row 44
$ws.Range("D44").Value = "'0.005314"
$ws.Range("G44").Value = "'9"
# Row 45
$ws.Range("D45").Value = "'0.00005236"
$ws.Range("G45").Value = "'9"
# Row 46
$ws.Range("G46").Value = "'9"
# Row 47
$ws.Range("D47").Value = "'0.6703"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOINBestin24h"
$ws.Range("G47").Value = "'9"
# Row 48
$ws.Range("D48").Value = "'0.02762"
$ws.Range("G48").Value = "'9"
# Row 49
$ws.Range("G49").Value = "'9"
# Row 50
$ws.Range("G50").Value = "'9"
# Row 51
$ws.Range("G51").Value = "'9"
